$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.719.21'
$ws.Cells.Item(2, 5).Value = '  +0.90%  '

$ws.Cells.Item(3, 4).Value = '1.889.41'
$ws.Cells.Item(3, 5).Value = '  +0.92%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.000'
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '248.46'
$ws.Cells.Item(5, 5).Value = '  +0.64%  '

$ws.Cells.Item(6, 5).Value = '  -0.05%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4739'
$ws.Cells.Item(7, 5).Value = '  +0.04%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2930'

$ws.Cells.Item(9, 5).Value = '  +0.45%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '22.02'
$ws.Cells.Item(10, 5).Value = '  +0.34%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07812'
$ws.Cells.Item(11, 5).Value = '  +1.22%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '97.01'
$ws.Cells.Item(12, 5).Value = '  -0.64%  '

$ws.Cells.Item(13, 4).Value = '1.891.60'
$ws.Cells.Item(13, 5).Value = '  +0.90%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.7351'
$ws.Cells.Item(14, 5).Value = '  -0.64%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '5.248'
$ws.Cells.Item(15, 5).Value = '  +2.45%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '284.63'
$ws.Cells.Item(16, 5).Value = '  +4.13%  '

$ws.Cells.Item(17, 4).Value = '30.950.71'
$ws.Cells.Item(17, 5).Value = '  +1.73%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '13.22'
$ws.Cells.Item(18, 5).Value = '  -1.36%  '

$ws.Cells.Item(19, 2).Value = 'Dai'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.001'
$ws.Cells.Item(19, 5).Value = '  -0.12%  '

$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.000007530'
$ws.Cells.Item(20, 5).Value = '  -0.09%  '

$ws.Cells.Item(21, 4).Value = '2.139.49'
$ws.Cells.Item(21, 5).Value = '  +0.56%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.332'
$ws.Cells.Item(22, 5).Value = '  +2.02%  '

$ws.Cells.Item(23, 5).Value = '  +0.01%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.256'
$ws.Cells.Item(24, 5).Value = '  +1.31%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.230'
$ws.Cells.Item(25, 5).Value = '  -0.26%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '164.29'
$ws.Cells.Item(26, 5).Value = '  +0.24%  '

$ws.Cells.Item(27, 5).Value = '  +0.54%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.926'
$ws.Cells.Item(28, 5).Value = '  -0.29%  '

$ws.Cells.Item(29, 5).Value = '  -1.70%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.09737'
$ws.Cells.Item(30, 5).Value = '  -3.26%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.498'
$ws.Cells.Item(31, 5).Value = '  -0.85%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.299'
$ws.Cells.Item(32, 5).Value = '  -0.31%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.189'
$ws.Cells.Item(33, 5).Value = '  +2.08%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.04864'
$ws.Cells.Item(34, 5).Value = '  +0.99%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.126'
$ws.Cells.Item(35, 5).Value = '  +0.23%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.6981'
$ws.Cells.Item(36, 5).Value = '  -0.02%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.724'
$ws.Cells.Item(37, 5).Value = '  +0.13%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01906'
$ws.Cells.Item(38, 5).Value = '  +2.48%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.802'
$ws.Cells.Item(39, 5).Value = '  +2.24%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.399'
$ws.Cells.Item(40, 5).Value = '  +1.78%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '76.17'
$ws.Cells.Item(41, 5).Value = '  +7.32%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.015'
$ws.Cells.Item(42, 5).Value = '  +1.87%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.4261'
$ws.Cells.Item(43, 5).Value = '  +1.61%  '

$ws.Cells.Item(44, 5).Value = '  -0.01%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.8337'
$ws.Cells.Item(45, 5).Value = '  -0.43%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '101.59'
$ws.Cells.Item(46, 5).Value = '  -1.07%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '9.455'
$ws.Cells.Item(47, 5).Value = '  +1.74%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '35.66'
$ws.Cells.Item(48, 5).Value = '  +0.35%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '7.037'
$ws.Cells.Item(49, 5).Value = '  +0.29%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '919.25'
$ws.Cells.Item(50, 5).Value = '  -0.05%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.05751'
